# Generate Report for Handoff
# Adds two new file rows (10d98b70-... and 27d6d7bd-...) to the
# Overview / zh-cn / de-de sheets, mirroring the existing "ht" rows.
#
# NOTE: every literal is written with a leading "'" (Excel's
# force-text prefix). Without it, values like "True"/"False" get
# auto-coerced to real booleans, and an empty "" assignment is
# dropped instead of producing an empty shared-string cell.
function Text([string]$s) { return "'" + $s }

$wb = $excel.ActiveWorkbook

$commitSha = "343350f330ca9a92de2af2f340fe01d78b02fcbf"
$repoBase  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/"

$files = @(
    @{ Name = "10d98b70-0f89-4fb0-affb-d3ac28d0d60c.md"; ZhXlf = "10d98b70-0f89-4fb0-affb-d3ac28d0d60c.784f0597682c7a2103cd59c9aaaa48f79e1ae6c9.zh-cn.xlf"; DeXlf = "10d98b70-0f89-4fb0-affb-d3ac28d0d60c.784f0597682c7a2103cd59c9aaaa48f79e1ae6c9.de-de.xlf" },
    @{ Name = "27d6d7bd-0eca-4740-a37e-fc21c2e6faf2.md"; ZhXlf = "27d6d7bd-0eca-4740-a37e-fc21c2e6faf2.39db929c5f1e87d4e75f767f18fba5f5167cf956.zh-cn.xlf"; DeXlf = "27d6d7bd-0eca-4740-a37e-fc21c2e6faf2.39db929c5f1e87d4e75f767f18fba5f5167cf956.de-de.xlf" }
)

$zhHandoffDateTime = "2016-11-03 18:59:56"
$deHandoffDateTime = "2016-11-03 19:00:10"
$overviewDateTime  = "2016-11-03 19:00:10"
$epoch             = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Sheet "Overview" -- File Name / Path And Name / Extension /
#   Publish URL / zh-cn / de-de / Latest HO Xliff Generate Date
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

foreach ($f in $files) {
    $null = $loOverview.ListRows.Add()
    $r = $loOverview.Range.Rows.Count + $loOverview.Range.Row - 1

    $wsOverview.Cells.Item($r, 1).Value = Text $f.Name
    $wsOverview.Cells.Item($r, 3).Value = Text ".md"
    $wsOverview.Cells.Item($r, 4).Value = Text ""
    $wsOverview.Cells.Item($r, 5).Value = Text "Ready for handoff"
    $wsOverview.Cells.Item($r, 6).Value = Text "Ready for handoff"
    $wsOverview.Cells.Item($r, 7).Value = Text $overviewDateTime
    $wsOverview.Cells.Item($r, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $null = $wsOverview.Hyperlinks.Add(
        $wsOverview.Cells.Item($r, 2),
        "$repoBase$($f.Name)",
        [Type]::Missing,
        [Type]::Missing,
        "e2e\$($f.Name)")
}

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

foreach ($f in $files) {
    $null = $loZh.ListRows.Add()
    $r = $loZh.Range.Rows.Count + $loZh.Range.Row - 1

    $wsZh.Cells.Item($r, 2).Value = Text ".md"
    $wsZh.Cells.Item($r, 3).Value = Text "Ready for handoff"
    $wsZh.Cells.Item($r, 4).Value = Text "e2e"
    $wsZh.Cells.Item($r, 5).Value = Text "ht"
    $wsZh.Cells.Item($r, 6).Value = Text "False"
    $wsZh.Cells.Item($r, 7).Value = Text $f.ZhXlf
    $wsZh.Cells.Item($r, 8).Value = Text $zhHandoffDateTime
    $wsZh.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $wsZh.Cells.Item($r, 9).Value = Text ""
    $wsZh.Cells.Item($r, 10).Value = Text ""
    $wsZh.Cells.Item($r, 11).Value = Text $epoch
    $wsZh.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $wsZh.Cells.Item($r, 12).Value = Text ""
    $wsZh.Cells.Item($r, 13).Value = Text "True"
    $wsZh.Cells.Item($r, 14).Value = Text ""
    $wsZh.Cells.Item($r, 15).Value = Text "False"
    $wsZh.Cells.Item($r, 16).Value = Text ""

    $null = $wsZh.Hyperlinks.Add(
        $wsZh.Cells.Item($r, 1),
        "$repoBase$($f.Name)",
        [Type]::Missing,
        [Type]::Missing,
        $f.Name)
}

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

foreach ($f in $files) {
    $null = $loDe.ListRows.Add()
    $r = $loDe.Range.Rows.Count + $loDe.Range.Row - 1

    $wsDe.Cells.Item($r, 2).Value = Text ".md"
    $wsDe.Cells.Item($r, 3).Value = Text "Ready for handoff"
    $wsDe.Cells.Item($r, 4).Value = Text "e2e"
    $wsDe.Cells.Item($r, 5).Value = Text "ht"
    $wsDe.Cells.Item($r, 6).Value = Text "False"
    $wsDe.Cells.Item($r, 7).Value = Text $f.DeXlf
    $wsDe.Cells.Item($r, 8).Value = Text $deHandoffDateTime
    $wsDe.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $wsDe.Cells.Item($r, 9).Value = Text ""
    $wsDe.Cells.Item($r, 10).Value = Text ""
    $wsDe.Cells.Item($r, 11).Value = Text $epoch
    $wsDe.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $wsDe.Cells.Item($r, 12).Value = Text ""
    $wsDe.Cells.Item($r, 13).Value = Text "True"
    $wsDe.Cells.Item($r, 14).Value = Text ""
    $wsDe.Cells.Item($r, 15).Value = Text "False"
    $wsDe.Cells.Item($r, 16).Value = Text ""

    $null = $wsDe.Hyperlinks.Add(
        $wsDe.Cells.Item($r, 1),
        "$repoBase$($f.Name)",
        [Type]::Missing,
        [Type]::Missing,
        $f.Name)
}

# ---------------------------------------------------------------
# Column-width touch-up (matches the widened "Status" / zh-cn / de-de
# columns once "Ready for handoff" is the longest value in them).
# The COM ColumnWidth setter snaps to the workbook's default-font
# character grid, so this lands on the closest attainable width.
# ---------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.9
$wsOverview.Columns.Item(6).ColumnWidth = 16.9
$wsZh.Columns.Item(3).ColumnWidth = 16.9
$wsDe.Columns.Item(3).ColumnWidth = 16.9

Write-Host "Handoff rows added."
